# Apply updates to column F (dSF) values as per repull/push of data + mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new value for column F
$updates = @{
    2  = -3
    3  = -5
    6  = -6
    7  = 2
    9  = -6
    10 = 2
    12 = -2
    13 = -5
    14 = -10
    16 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
